$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record needs to be inserted above row 555, pushing the
# existing rows (555-581) down by one (556-582).
$ws.Rows.Item(555).EntireRow.Insert()

# The inserted row is blank; clone the row that is now directly below it
# (old row 555, now at 556) so every column keeps the same formatting /
# static values (Mercado, Region, Categoria, etc.), then overwrite the
# cells that actually hold new data for this record.
$ws.Range("A556:R556").Copy()
$ws.Range("A555:R555").PasteSpecial()

$ws.Range("D555").Value = 45147
$ws.Range("J555").Value = 560
$ws.Range("K555").Value = 5500
$ws.Range("L555").Value = 6000
$ws.Range("M555").Value = 5750
$ws.Range("P555").Value = 288
